$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.341.91"
$ws.Range("E2").Value = "  +0.77%  "

$ws.Range("D3").Value = "2.187.88"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.626"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "68.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.43%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.97%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0937"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.69%  "

$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("D15").Value = "2.514.93"
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.873"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.26%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.37%  "

$ws.Range("D18").Value = "2.196.70"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "41.274.30"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").Value = "0.0₃0954"
$ws.Range("E20").Value = "  +1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("E24").Value = "  +3.48%  "

$ws.Range("B25").Value = "WEMIXToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.01%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +22.64%  "

$ws.Range("E28").Value = "  +5.47%  "

$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.119"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0743"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.38%  "

$ws.Range("E34").Value = "  +0.47%  "

$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "27.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +17.58%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.17%  "

$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.96%  "

$ws.Range("E39").Value = "  +13.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +27.31%  "

$ws.Range("E41").Value = "  -0.99%  "

$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.12%  "

$ws.Range("E45").Value = "  +5.92%  "

$ws.Range("E46").Value = "  +0.19%  "

$ws.Range("E47").Value = "  +3.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("E49").Value = "  +5.45%  "

$ws.Range("E50").Value = "  +1.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.20%  "
